# This workbook holds loan-product scenario data on "ProductLoanInput" and
# "ProductLoanOutput". We re-assign two field values on ProductLoanInput and
# leave ProductLoanInput as the active/selected sheet (previously
# ProductLoanOutput was active).

$wb = $excel.ActiveWorkbook
$wsIn = $wb.Worksheets.Item("ProductLoanInput")

# repaymentstrategy (row 17) used to be "RBI (India)"; it now gets a brand
# new value, rendered left/top aligned (a new cell style).
$rB17 = $wsIn.Range("B17")
$rB17.Value = "Penalties, Fees, Interest, Principal order"
$rB17.HorizontalAlignment = -4131   # xlLeft
$rB17.VerticalAlignment = -4160     # xlTop

# preclosureinterestcalculationrule (row 22) now takes the value that used
# to belong to repaymentstrategy.
$wsIn.Range("B22").Value = "RBI (India)"

# ProductLoanInput becomes the active sheet, with B17 selected.
$wsIn.Activate()
$wsIn.Range("B17").Select()
